$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 450.5
$ws.Range("I107").Value = 455.64285
$ws.Range("K107").Value = 455.64285
$ws.Range("M107").Value = 1464.35715
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H132").Value = 112898.78
$ws.Range("I132").Value = 68425.60000000001
$ws.Range("K132").Value = 205276.8
$ws.Range("M132").Value = -202746.8
$ws.Range("H138").Value = 3871.1526
$ws.Range("J138").Value = 4623.5776
$ws.Range("L138").Value = 13870.7328
$ws.Range("N138").Value = -24150.7328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 2055
$ws.Range("J29").Value = 2055
$ws.Range("L29").Value = 2055
$ws.Range("N29").Value = -2671
$ws.Range("H53").Value = 16995
$ws.Range("I53").Value = 16995
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 16995
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -16313
$ws.Range("N53").ClearContents()
$ws.Range("H131").Value = 75712.86
$ws.Range("J131").Value = 75712.86
$ws.Range("L131").Value = 75712.86
$ws.Range("N131").Value = -85792.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3320.762
$ws.Range("I20").Value = 2730.6667
$ws.Range("J20").Value = 4107.5557
$ws.Range("K20").Value = 2730.6667
$ws.Range("L20").Value = 4107.5557
$ws.Range("M20").Value = -2483.6667
$ws.Range("N20").Value = -4601.5557
$ws.Range("H63").Value = 50250
$ws.Range("J63").Value = 50250
$ws.Range("L63").Value = 50250
$ws.Range("N63").Value = -51622
$ws.Range("H66").Value = 50250
$ws.Range("J66").Value = 50250
$ws.Range("L66").Value = 150750
$ws.Range("N66").Value = -157614
$ws.Range("H99").Value = 6103.6665
$ws.Range("I99").Value = 1942.2
$ws.Range("J99").Value = 11305.5
$ws.Range("K99").Value = 1942.2
$ws.Range("L99").Value = 11305.5
$ws.Range("M99").Value = -444.2
$ws.Range("N99").Value = -14301.5
$ws.Range("H134").Value = 2329.5898
$ws.Range("I134").Value = 2072.9714
$ws.Range("K134").Value = 6218.914199999999
$ws.Range("M134").Value = -3683.914199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 26990
$ws.Range("J57").Value = 26990
$ws.Range("L57").Value = 26990
$ws.Range("N57").Value = -28110
$ws.Range("H105").Value = 1833
$ws.Range("I105").Value = 1874.5
$ws.Range("J105").Value = 1750
$ws.Range("K105").Value = 1874.5
$ws.Range("L105").Value = 1750
$ws.Range("M105").Value = -127.5
$ws.Range("N105").Value = -5244
$ws.Range("H132").Value = 4929.875
$ws.Range("J132").Value = 5056.8
$ws.Range("L132").Value = 15170.4
$ws.Range("N132").Value = -20230.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4392.8945
$ws.Range("I5").Value = 605.625
$ws.Range("J5").Value = 7147.273
$ws.Range("K5").Value = 1816.875
$ws.Range("L5").Value = 21441.819
$ws.Range("M5").Value = -1704.875
$ws.Range("N5").Value = -21665.819
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372
$ws.Range("H63").Value = 1200
$ws.Range("I63").Value = 1200
$ws.Range("K63").Value = 3600
$ws.Range("M63").Value = -2851
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 2500
$ws.Range("K64").Value = 7500
$ws.Range("M64").Value = -7230
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864
$ws.Range("H66").Value = 1200
$ws.Range("I66").Value = 1200
$ws.Range("K66").Value = 10800
$ws.Range("M66").Value = -7056
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 2500
$ws.Range("K67").Value = 7500
$ws.Range("M67").Value = -6564
$ws.Range("H69").Value = 4211.25
$ws.Range("I69").Value = 3948.3333
$ws.Range("K69").Value = 11844.9999
$ws.Range("M69").Value = -11033.9999
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5685
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 4211.25
$ws.Range("I72").Value = 3948.3333
$ws.Range("K72").Value = 35534.9997
$ws.Range("M72").Value = -31478.9997
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4908
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 1500
$ws.Range("I74").Value = 1500
$ws.Range("K74").Value = 4500
$ws.Range("M74").Value = -3439
$ws.Range("H75").Value = 15000
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61996
$ws.Range("H77").Value = 1500
$ws.Range("I77").Value = 1500
$ws.Range("K77").Value = 13500
$ws.Range("M77").Value = -8196
$ws.Range("H78").Value = 15000
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189984
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 90000
$ws.Range("N87").Value = -92496
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 270000
$ws.Range("N90").Value = -282480
$ws.Range("H107").Value = 1073.6
$ws.Range("I107").Value = 270
$ws.Range("J107").Value = 1274.5
$ws.Range("K107").Value = 810
$ws.Range("L107").Value = 3823.5
$ws.Range("M107").Value = 1110
$ws.Range("N107").Value = -7663.5
$ws.Range("H113").Value = 406.75
$ws.Range("J113").Value = 375.66666
$ws.Range("L113").Value = 1126.99998
$ws.Range("N113").Value = -5466.999980000001
$ws.Range("H135").Value = 4392.8945
$ws.Range("I135").Value = 605.625
$ws.Range("J135").Value = 7147.273
$ws.Range("K135").Value = 5450.625
$ws.Range("L135").Value = 64325.457
$ws.Range("M135").Value = -2915.625
$ws.Range("N135").Value = -69395.45699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4228.5
$ws.Range("I80").Value = 3891.25
$ws.Range("J80").Value = 4903
$ws.Range("K80").Value = 3891.25
$ws.Range("L80").Value = 4903
$ws.Range("M80").Value = -2893.25
$ws.Range("N80").Value = -6899
$ws.Range("H83").Value = 4228.5
$ws.Range("I83").Value = 3891.25
$ws.Range("J83").Value = 4903
$ws.Range("K83").Value = 19456.25
$ws.Range("L83").Value = 24515
$ws.Range("M83").Value = -14464.25
$ws.Range("N83").Value = -34499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 75008730
$ws.Range("J18").Value = 75008730
$ws.Range("L18").Value = 75008730
$ws.Range("N18").Value = -75009076
$ws.Range("H132").Value = 2881.6206
$ws.Range("J132").Value = 2503.375
$ws.Range("L132").Value = 7510.125
$ws.Range("N132").Value = -12570.125
